# Normalize the "Recorded By" (column G) entries so that the two-name
# combinations involving "dnasr281@gmail.com" and "admin@admin.com" have
# their display order swapped (e.g. "System, dnasr281@gmail.com" becomes
# "dnasr281@gmail.com, System"). Other combinations - such as the
# "backup@backdoor.com, System" pairing or single-name entries - are left
# untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$targets = @(
    "System, dnasr281@gmail.com",
    "admin@admin.com, System",
    "admin@admin.com, dnasr281@gmail.com"
)

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)
    $value = $cell.Value2

    if ($null -eq $value) {
        continue
    }

    if ($targets -contains $value) {
        $parts = $value -split ', '
        $reversed = ($parts[($parts.Length - 1)..0]) -join ', '
        $cell.Value2 = $reversed
    }
}
